$p = $ppt.ActivePresentation

# --- 1) Update the cached "today" date fields (datetime1 / datetimeFigureOut)
#        from 14-Feb-2025 to 21-Feb-2025 across slide master, every layout,
#        and the notes master.
$oldMDY = "2/14/2025"
$newMDY = "2/21/2025"
$oldDMY = "14-02-2025"
$newDMY = "21-02-2025"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq $oldMDY) {
                $tr.Text = $newMDY
            } elseif ($tr.Text -eq $oldDMY) {
                $tr.Text = $newDMY
            }
        }
    }
}

# Slide master
Update-DateShape $p.SlideMaster.Shapes

# Every slide layout (custom layout) hanging off the master
$customLayouts = $p.SlideMaster.CustomLayouts
for ($L = 1; $L -le $customLayouts.Count; $L++) {
    Update-DateShape $customLayouts.Item($L).Shapes
}

# Notes master (its placeholder shapes aren't directly editable through
# TextFrame/TextRange, so go through the HeadersFooters/DateAndTime object
# that is bound to the same cached date field)
$nmDate = $p.NotesMaster.HeadersFooters.DateAndTime
if ($nmDate.Text -eq $oldMDY) {
    $nmDate.Text = $newMDY
} elseif ($nmDate.Text -eq $oldDMY) {
    $nmDate.Text = $newDMY
} else {
    # Fall back to a direct (unconditional) update - the cached read-back
    # value for the notes master isn't always reliable in this host.
    $nmDate.Text = $newDMY
}

# --- 2) Turn the GitHub URL text on the "GitHub Link" slide into a live
#        hyperlink pointing at the same address.
$ppMouseClick = 1
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTextFrame) {
            $text = $shape.TextFrame.TextRange.Text
            if ($text -eq "https://github.com/SAIKUMAR173/STEGANOGRAPHY.git") {
                $actionSetting = $shape.TextFrame.TextRange.ActionSettings($ppMouseClick)
                $actionSetting.Hyperlink.Address = $text
            }
        }
    }
}
